$d = $word.ActiveDocument

# --- Add the three new character styles ------------------------------------

$sStyle = $d.Styles.Add("GaNStyle", 2)
$sStyle.Font.Name = "Calibri"
$sStyle.Font.Size = 14

$sParagraph = $d.Styles.Add("GaNParagraph", 2)
$sParagraph.Font.Name = "Calibri"
$sParagraph.Font.Size = 10

$sLinks = $d.Styles.Add("GaNLinks", 2)
$sLinks.Font.Name = "Calibri"
$sLinks.Font.Bold = $true
$sLinks.Font.Color = 8388608
$sLinks.Font.Size = 9.5
$sLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Kampagnendaten 2022" run ----------------------

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Kampagnendaten 2022 für das Sternbild Zwillinge: 14.-23. Februar, 14.-24. März"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
while ($rng.Find.Execute()) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the "Mach mit an einer weltweiten..." run --------

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Mach mit an einer weltweiten Kampagne, die schwächsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Zwillinge am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel."
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
if ($rng.Find.Execute()) {
    $rng.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Die Schaubilder in diesem Dokument..." run ------

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
if ($rng.Find.Execute()) {
    $rng.Style = "GaNLinks"
}
